# [FIX] filter charge_type (expense, gl, revenue and trial balance web)
#
# Adds a new "Charge Type" filter-criteria row to the report header block,
# just above the existing "Run By" / "Run Date" rows, so the Expense Ledger
# report can be filtered by charge type.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current "Run By" row (row 10). This shifts the
# existing "Run By" (10->11), "Run Date" (11->12), the blank spacer row
# (12->13) and the detail-table header row (13->14) down by one, carrying
# their original formatting with them.
$ws.Rows.Item(10).Insert()

# The new row 10 inherited row 9's column layout (a two-field date row), so
# restore the single-label/single-input layout used by the other filter rows:
#   Column A: bold label cell (same style as "Account Code", "Budget", ...)
#   Column B: the (empty) input cell next to the label
#   Columns C.. : plain/unformatted cells, matching every other filter row.

# A10 <- label style copied from another one-field filter label (A3 "Account Code")
$ws.Range("A3").Copy()
$ws.Range("A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("A10").Value = "Charge Type"

# B10 <- input-cell style copied from an existing filter input cell (B9 "Date To" input)
$ws.Range("B9").Copy()
$ws.Range("B10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# C10:AX10 <- plain/default style, matching the rest of the row's trailing cells
$ws.Range("C1").Copy()
$ws.Range("C10:AX10").PasteSpecial(-4122)
$excel.CutCopyMode = $false
